$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue 'D2' '28.411.75'
Set-TextValue 'E2' '  +3.42%  '
Set-TextValue 'D3' '1.867.68'
Set-TextValue 'E3' '  +1.91%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '338.66'
Set-TextValue 'E5' '  +2.19%  '
Set-TextValue 'D6' '1.001'
Set-TextValue 'E6' '  -0.02%  '
Set-TextValue 'D7' '0.4688'
Set-TextValue 'E7' '  +1.43%  '
Set-TextValue 'D8' '0.3962'
Set-TextValue 'E8' '  +3.29%  '
Set-TextValue 'D9' '47.71'
Set-TextValue 'E9' '  +2.23%  '
Set-TextValue 'D10' '0.08018'
Set-TextValue 'E10' '  +1.61%  '
Set-TextValue 'D11' '0.9992'
Set-TextValue 'E11' '  +2.75%  '
Set-TextValue 'D12' '21.99'
Set-TextValue 'E12' '  +4.04%  '
Set-TextValue 'D13' '6.028'
Set-TextValue 'E13' '  +2.36%  '
Set-TextValue 'D14' '1.865.01'
Set-TextValue 'E14' '  +1.60%  '
Set-TextValue 'D15' '7.260'
Set-TextValue 'E15' '  +2.79%  '
Set-TextValue 'D16' '90.70'
Set-TextValue 'E16' '  +3.00%  '
Set-TextValue 'E17' '  +0.06%  '
Set-TextValue 'D18' '0.00001043'
Set-TextValue 'E18' '  +1.33%  '
Set-TextValue 'D19' '0.06621'
Set-TextValue 'E19' '  +0.00%  '
Set-TextValue 'D20' '17.52'
Set-TextValue 'E20' '  +2.86%  '
Set-TextValue 'D21' '1.001'
Set-TextValue 'E21' '  -0.08%  '
Set-TextValue 'D22' '28.424.57'
Set-TextValue 'E22' '  +3.46%  '
Set-TextValue 'D23' '5.465'
Set-TextValue 'E23' '  +2.24%  '
Set-TextValue 'D24' '11.04'
Set-TextValue 'E24' '  +2.07%  '
Set-TextValue 'D25' '2.264'
Set-TextValue 'E25' '  -2.00%  '
Set-TextValue 'D26' '2.081.75'
Set-TextValue 'E26' '  +1.31%  '
Set-TextValue 'D27' '160.71'
Set-TextValue 'E27' '  +2.24%  '
Set-TextValue 'D28' '19.79'
Set-TextValue 'E28' '  +2.06%  '
Set-TextValue 'D29' '2.117'
Set-TextValue 'E29' '  +2.46%  '
Set-TextValue 'D30' '5.495'
Set-TextValue 'E30' '  +3.59%  '
Set-TextValue 'D31' '120.14'
Set-TextValue 'E31' '  +0.89%  '
Set-TextValue 'D32' '0.9847'
Set-TextValue 'E32' '  +2.82%  '
Set-TextValue 'D33' '0.09492'
Set-TextValue 'E33' '  +2.15%  '
Set-TextValue 'D34' '3.588'
Set-TextValue 'E34' '  +0.60%  '
Set-TextValue 'B35' 'Filecoin'
Set-TextValue 'C35' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D35' '5.353'
Set-TextValue 'E35' '  +2.16%  '
Set-TextValue 'B36' 'ARBITRUM'
Set-TextValue 'C36' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D36' '1.373'
Set-TextValue 'E36' '  +4.28%  '
Set-TextValue 'D37' '0.06107'
Set-TextValue 'E37' '  +2.81%  '
Set-TextValue 'D38' '0.02247'
Set-TextValue 'E38' '  +2.25%  '
Set-TextValue 'D39' '8.328'
Set-TextValue 'E39' '  +3.06%  '
Set-TextValue 'D40' '1.178'
Set-TextValue 'E40' '  +1.54%  '
Set-TextValue 'D41' '0.5940'
Set-TextValue 'E41' '  +2.32%  '
Set-TextValue 'D42' '1.001'
Set-TextValue 'E42' '  +0.00%  '
Set-TextValue 'D43' '0.1875'
Set-TextValue 'E43' '  +1.85%  '
Set-TextValue 'D44' '10.33'
Set-TextValue 'E44' '  +2.93%  '
Set-TextValue 'D45' '1.282'
Set-TextValue 'E45' '  -0.30%  '
Set-TextValue 'D46' '0.5572'
Set-TextValue 'E46' '  +1.47%  '
Set-TextValue 'D47' '12.16'
Set-TextValue 'E47' '  +1.45%  '
Set-TextValue 'D48' '1.955'
Set-TextValue 'E48' '  +4.31%  '
Set-TextValue 'D49' '0.06983'
Set-TextValue 'E49' '  +5.01%  '
Set-TextValue 'D50' '2.072'
Set-TextValue 'E50' '  +14.43%  '
Set-TextValue 'D51' '111.76'
Set-TextValue 'E51' '  +1.12%  '
